$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-14 07:04:20"

# --- Worksheet index 2 ---
$ws = $wb.Worksheets.Item(2)

$ws.Range("AA2").Value = $newTimestamp

$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 427
$ws.Range("E3").Value = 189
$ws.Range("F3").Value = 238
$ws.Range("G3").Value = 17.79
$ws.Range("H3").Value = 7.88
$ws.Range("I3").Value = 9.92
$ws.Range("J3").Value = 92
$ws.Range("K3").Value = 99
$ws.Range("V3").Value = 14
$ws.Range("AA3").Value = $newTimestamp

$ws.Range("AA4").Value = $newTimestamp

$ws.Range("AA5").Value = $newTimestamp

$ws.Range("AA6").Value = $newTimestamp

$ws.Range("C7").Value = 16
$ws.Range("D7").Value = 211
$ws.Range("E7").Value = 94
$ws.Range("F7").Value = 117
$ws.Range("G7").Value = 13.19
$ws.Range("H7").Value = 5.88
$ws.Range("I7").Value = 7.31
$ws.Range("J7").Value = 47
$ws.Range("K7").Value = 41
$ws.Range("AA7").Value = $newTimestamp

$ws.Range("AA8").Value = $newTimestamp

$ws.Range("C9").Value = 25
$ws.Range("D9").Value = 384
$ws.Range("E9").Value = 207
$ws.Range("F9").Value = 177
$ws.Range("G9").Value = 15.36
$ws.Range("H9").Value = 8.28
$ws.Range("I9").Value = 7.08
$ws.Range("J9").Value = 101
$ws.Range("K9").Value = 86
$ws.Range("V9").Value = 12
$ws.Range("AA9").Value = $newTimestamp

$ws.Range("AA10").Value = $newTimestamp

$ws.Range("AA11").Value = $newTimestamp

$ws.Range("AA12").Value = $newTimestamp

$ws.Range("AA13").Value = $newTimestamp

$ws.Range("AA14").Value = $newTimestamp

$ws.Range("AA15").Value = $newTimestamp

$ws.Range("AA16").Value = $newTimestamp

$ws.Range("AA17").Value = $newTimestamp

$ws.Range("AA18").Value = $newTimestamp

$ws.Range("AA19").Value = $newTimestamp

$ws.Range("AA20").Value = $newTimestamp

$ws.Range("AA21").Value = $newTimestamp

$ws.Range("AA22").Value = $newTimestamp

$ws.Range("AA23").Value = $newTimestamp

$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 409
$ws.Range("E24").Value = 181
$ws.Range("F24").Value = 228
$ws.Range("G24").Value = 16.36
$ws.Range("H24").Value = 7.24
$ws.Range("I24").Value = 9.12
$ws.Range("J24").Value = 88
$ws.Range("K24").Value = 99
$ws.Range("AA24").Value = $newTimestamp

$ws.Range("AA25").Value = $newTimestamp

$ws.Range("AA26").Value = $newTimestamp

# --- Worksheet index 3 ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("AA2").Value = $newTimestamp

$ws.Range("AA3").Value = $newTimestamp

$ws.Range("AA4").Value = $newTimestamp

$ws.Range("C5").Value = 11
$ws.Range("D5").Value = 158
$ws.Range("E5").Value = 84
$ws.Range("F5").Value = 74
$ws.Range("G5").Value = 14.36
$ws.Range("H5").Value = 7.64
$ws.Range("I5").Value = 6.73
$ws.Range("J5").Value = 42
$ws.Range("K5").Value = 37
$ws.Range("V5").Value = 10
$ws.Range("AA5").Value = $newTimestamp

$ws.Range("AA6").Value = $newTimestamp

$ws.Range("AA7").Value = $newTimestamp

$ws.Range("AA8").Value = $newTimestamp

$ws.Range("C9").Value = 23
$ws.Range("D9").Value = 440
$ws.Range("E9").Value = 189
$ws.Range("F9").Value = 251
$ws.Range("G9").Value = 19.13
$ws.Range("H9").Value = 8.22
$ws.Range("I9").Value = 10.91
$ws.Range("J9").Value = 82
$ws.Range("K9").Value = 103
$ws.Range("AA9").Value = $newTimestamp

$ws.Range("AA10").Value = $newTimestamp

$ws.Range("AA11").Value = $newTimestamp

$ws.Range("AA12").Value = $newTimestamp

$ws.Range("AA13").Value = $newTimestamp

$ws.Range("AA14").Value = $newTimestamp

$ws.Range("C15").Value = 20
$ws.Range("D15").Value = 401
$ws.Range("E15").Value = 211
$ws.Range("F15").Value = 190
$ws.Range("G15").Value = 20.05
$ws.Range("H15").Value = 10.55
$ws.Range("I15").Value = 9.5
$ws.Range("J15").Value = 83
$ws.Range("K15").Value = 75
$ws.Range("V15").Value = 8
$ws.Range("AA15").Value = $newTimestamp

$ws.Range("AA16").Value = $newTimestamp

$ws.Range("AA17").Value = $newTimestamp

$ws.Range("AA18").Value = $newTimestamp

$ws.Range("AA19").Value = $newTimestamp

$ws.Range("AA20").Value = $newTimestamp

$ws.Range("AA21").Value = $newTimestamp

$ws.Range("AA22").Value = $newTimestamp

$ws.Range("AA23").Value = $newTimestamp

$ws.Range("AA24").Value = $newTimestamp

$ws.Range("AA25").Value = $newTimestamp

$ws.Range("AA26").Value = $newTimestamp
